$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '92.956.62'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.66%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.109.32'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.57%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.24'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.19%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '612.05'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.03%  '

$ws.Range('E7').Value = '  -0.33%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.394'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +5.46%  '

$ws.Range('E9').Value = '  -0.04%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.108.66'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +30.13%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.752'
$ws.Range('D11').ClearFormats()

$ws.Range('E12').Value = '  -0.86%  '

$ws.Range('E13').Value = '  +1.80%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '93.018.26'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.98%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.33'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.46%  '

$ws.Range('E16').Value = '  -1.49%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.704.34'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.97%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.115.23'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.26%  '

$ws.Range('E19').Value = '  +0.98%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.78'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.27%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.77'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.31%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '446.49'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.07%  '

$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.28'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.75%  '

$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000201'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.22%  '

$ws.Range('E25').Value = '  -3.28%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '86.85'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.33%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.68'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.89%  '

$ws.Range('E29').Value = '  +0.00%  '

$ws.Range('E30').Value = '  +10.78%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.232'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.21%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.170'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.66%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.15'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.88%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '8.05'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.96%  '

$ws.Range('E35').Value = '  -6.68%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '26.03'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.06%  '

$ws.Range('E37').Value = '  -6.32%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '490.77'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.99%  '

$ws.Range('E39').Value = '  -1.76%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.85'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.84%  '

$ws.Range('E41').Value = '  -4.75%  '

$ws.Range('E42').Value = '  -3.54%  '

$ws.Range('E43').Value = '  +4.27%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.37'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.84%  '

$ws.Range('E45').Value = '  +0.01%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '163.44'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.37%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.90'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.31%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.685'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.11%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.38'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.68%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0332'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +4.57%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '44.04'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.11%  '
